# Insert 5 new rows before row 429 with a new weekly price report for
# "Provincia del Elquí" (Fruta / Chirimoya / Cultivar IV Región), pushing
# the previously-existing rows 429-460 down to 434-465.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows at row 429 (existing rows 429+ shift down).
$insertRange = $ws.Range("A429:T433")
$insertRange.EntireRow.Insert()

# Constant columns (same for every Chirimoya / Cultivar IV Región row).
$marketId = 9
$market = "Vega Central Mapocho de Santiago"
$region = "Metropolitana"
$codreg = 13
$tipo = "Fruta"
$productoId = 100107
$producto = "Otros"
$categoriaId = 100107002
$categoria = "Chirimoya"
$variedad = "Cultivar IV Región"

# Data for the new weekly snapshot (date 45267, Provincia del Elquí,
# $/bandeja 10 kilos).
$newRows = @(
    @{ Row = 429; Calidad = "Especial";                Volumen = 90;  PMin = 21000; PMax = 21000; PProm = 21000; Unidad = "$/bandeja 10 kilos"; Origen = "Provincia del Elquí"; PKg = 2100; KgUnidad = 10 },
    @{ Row = 430; Calidad = "Extra (doble especial)";  Volumen = 90;  PMin = 24000; PMax = 24000; PProm = 24000; Unidad = "$/bandeja 10 kilos"; Origen = "Provincia del Elquí"; PKg = 2400; KgUnidad = 10 },
    @{ Row = 431; Calidad = "Primera";                 Volumen = 135; PMin = 18000; PMax = 18000; PProm = 18000; Unidad = "$/bandeja 10 kilos"; Origen = "Provincia del Elquí"; PKg = 1800; KgUnidad = 10 },
    @{ Row = 432; Calidad = "Segunda";                 Volumen = 90;  PMin = 14000; PMax = 14000; PProm = 14000; Unidad = "$/bandeja 10 kilos"; Origen = "Provincia del Elquí"; PKg = 1400; KgUnidad = 10 },
    @{ Row = 433; Calidad = "Tercera";                 Volumen = 45;  PMin = 12000; PMax = 12000; PProm = 12000; Unidad = "$/bandeja 10 kilos"; Origen = "Provincia del Elquí"; PKg = 1200; KgUnidad = 10 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $marketId
    $ws.Cells.Item($row, 2).Value = $market
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = 45267
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
